$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Drop the trailing "dist_to_shore / concavity / sd_SAL / log_sd_VEL" beta
#        columns (W:AB). This shifts nothing else since they were the last
#        columns on the sheet; it also prunes their now-unused shared strings.
$ws.Range("W1:AB5").EntireColumn.Delete()

# --- 2. Rename the remaining beta columns (now S:V) to the new variable set.
$ws.Range("S1").Value2 = "beta_mean_CHL"
$ws.Range("T1").Value2 = "beta_mean_SSH"
$ws.Range("U1").Value2 = "sd_beta_mean_CHL"
$ws.Range("V1").Value2 = "sd_beta_mean_SSH"

# --- 3. Rewrite the four data rows (model results) with the refreshed values.
$rows = @(
    @("without spatial",1.06,164.1,0.46,0.5,0.23,0.5,0.03,0.5,1248,316,726,2290,1247,320,729,2296,1,1.48,0.83,0.72,0.43),
    @("with spatial exp",2.09,42.6,0.48,0.51,0.25,0.51,0.03,0.5,1247,313,723,2283,1277,318,736,2331,2,1.43,0.96,0.7,0.54),
    @("with sp shpere",2.13,50.3,0.46,0.5,0.24,0.5,0.03,0.5,1248,315,726,2289,1277,318,735,2330,3,1.53,0.93,0.75,0.48),
    @("with sp gaussian",6.807,10.2,0.47,0.5,0.24,0.5,0.03,0.51,1247,315,725,2287,1276,318,737,2331,4,1.45,0.87,0.71,0.49)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $data = $rows[$r]
    $rowNum = $r + 2
    for ($c = 0; $c -lt $data.Length; $c++) {
        $colNum = $c + 1
        $ws.Cells.Item($rowNum, $colNum).Value2 = $data[$c]
    }
}

# --- 4. Update the conditional-formatting "below threshold" expression rules
#        to match the refreshed waic/CV columns (J:Q). FormatConditions.Type
#        2 == xlExpression (1 == xlCellValue, 3 == xlColorScale).
$thresholds = @{
    "J" = 1250
    "K" = 316
    "L" = 726
    "M" = 2286
    "N" = 1250
    "O" = 321
    "P" = 732
    "Q" = 2299
}

foreach ($col in $thresholds.Keys) {
    $rng = $ws.Range("$($col)2:$($col)5")
    foreach ($fc in $rng.FormatConditions) {
        if ($fc.Type -eq 2) {
            $fc.Formula1 = "$($col)2<$($thresholds[$col])"
        }
    }
}
